$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.431.22"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "1.842.80"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D5").Value = "'315.20"
$ws.Range("E5").Value = "  +1.83%  "
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("D7").Value = "'0.4775"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "'0.07470"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "'0.8874"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "'20.49"
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").Value = "1.844.84"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("D13").Value = "'0.07376"
$ws.Range("E13").Value = "  +4.25%  "
$ws.Range("D14").Value = "'5.489"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").Value = "'93.34"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "'6.609"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").Value = "'1.015"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("D21").Value = "27.443.76"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "'5.359"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "'10.75"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").Value = "2.090.85"
$ws.Range("E24").Value = "  +1.81%  "
$ws.Range("D25").Value = "'1.905"
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'152.69"
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("D28").Value = "'2.174"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").Value = "'0.08990"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "'0.7603"
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'4.569"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "'2.953"
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").Value = "'1.108"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("D38").Value = "'0.05374"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "'0.01970"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "'3.003"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D42").Value = "'0.5367"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "'2.383"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("D44").Value = "'0.1668"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "'8.573"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").Value = "'0.4988"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").Value = "'10.59"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "'1.015"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'105.32"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").Value = "'0.06327"
$ws.Range("E51").Value = "  +0.37%  "
